$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.162369132041931
$ws.Range("B1").Value = 2.369465589523315
$ws.Range("D1").Value = 2.395141124725342
$ws.Range("E1").Value = 1.218288779258728
